$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-11 with the new set of reddit post links (row 8 is unchanged).
$ws.Range("A2").Value = "/r/EarthPorn/comments/gz78ba/open_letter_to_steve_huffman_and_the_board_of/"
$ws.Range("A3").Value = "/r/EarthPorn/comments/hdsj7b/tried_to_find_a_leprechauns_pot_of_gold_but_got/"
$ws.Range("A4").Value = "/r/goodnews/comments/gwtp43/whats_new_content_creators_june_2020/"
$ws.Range("A5").Value = "/r/goodnews/comments/hc0d3l/good_news_its_friday_whats_your_feelgood_story/"
$ws.Range("A6").Value = "/r/learnpython/comments/hdhm9j/ask_anything_monday_weekly_thread/"
$ws.Range("A7").Value = "/r/learnpython/comments/hduh6k/how_to_approach_automate_the_boring_stuff/"
$ws.Range("A9").Value = "/r/pics/comments/hdyxvr/bubba_wallace_nascars_only_black_driver_with/"
$ws.Range("A10").Value = "/r/Python/comments/gdfaip/rpython_job_board_for_may_june_july/"
$ws.Range("A11").Value = "/r/Python/comments/ha4qfy/whats_everyone_working_on_this_week/"

# Remove the now-obsolete trailing rows (12-23) so the sheet shrinks back down.
$ws.Range("A12:A23").EntireRow.Delete()

# Restore the view: select column A from the top instead of the prior
# scrolled-down I50 selection.
$ws.Columns.Item(1).Select()
